$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '34.356.97'
$ws.Cells.Item(2, 5).Value = '  +0.53%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.791.75'
$ws.Cells.Item(3, 5).Value = '  +0.26%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '226.50'
$ws.Cells.Item(5, 5).Value = '  +0.09%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.555'
$ws.Cells.Item(6, 5).Value = '  +1.59%  '
$ws.Cells.Item(7, 5).Value = '  -0.04%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '32.35'
$ws.Cells.Item(8, 5).Value = '  +1.07%  '
$ws.Cells.Item(9, 5).Value = '  +1.18%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0690'
$ws.Cells.Item(10, 5).Value = '  -0.08%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '2.050.79'
$ws.Cells.Item(12, 5).Value = '  +0.27%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '11.00'
$ws.Cells.Item(13, 5).Value = '  -2.07%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '1.785.69'
$ws.Cells.Item(14, 5).Value = '  +0.14%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.628'
$ws.Cells.Item(15, 5).Value = '  +1.38%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '34.314.58'
$ws.Cells.Item(16, 5).Value = '  +0.58%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '4.21'
$ws.Cells.Item(17, 5).Value = '  +0.17%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '68.10'
$ws.Cells.Item(18, 5).Value = '  +0.13%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.0₃0800'
$ws.Cells.Item(19, 5).Value = '  +2.62%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '246.66'
$ws.Cells.Item(20, 5).Value = '  +0.03%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.94'
$ws.Cells.Item(21, 5).Value = '  +1.04%  '
$ws.Cells.Item(22, 5).Value = '  -0.10%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '4.17'
$ws.Cells.Item(23, 5).Value = '  +1.74%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.07'
$ws.Cells.Item(24, 5).Value = '  +0.87%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '162.21'
$ws.Cells.Item(25, 5).Value = '  +0.37%  '
$ws.Cells.Item(26, 5).Value = '  +0.30%  '
$ws.Cells.Item(27, 5).Value = '  +0.13%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.116'
$ws.Cells.Item(28, 5).Value = '  +1.64%  '
$ws.Cells.Item(29, 5).Value = '  +0.10%  '
$ws.Cells.Item(30, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '3.92'
$ws.Cells.Item(30, 5).Value = '  +8.58%  '
$ws.Cells.Item(31, 2).Value = 'PancakeSwap'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.23'
$ws.Cells.Item(31, 5).Value = '  -0.09%  '
$ws.Cells.Item(32, 2).Value = 'Hedera'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.0521'
$ws.Cells.Item(32, 5).Value = '  +0.37%  '
$ws.Cells.Item(33, 2).Value = 'Filecoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '3.79'
$ws.Cells.Item(33, 5).Value = '  +3.32%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.82'
$ws.Cells.Item(34, 5).Value = '  +0.12%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.441.94'
$ws.Cells.Item(35, 5).Value = '  -0.68%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.61'
$ws.Cells.Item(36, 5).Value = '  +9.21%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.662'
$ws.Cells.Item(37, 5).Value = '  +2.47%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.06'
$ws.Cells.Item(38, 5).Value = '  +1.64%  '
$ws.Cells.Item(39, 5).Value = '  -1.47%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '82.27'
$ws.Cells.Item(40, 5).Value = '  +2.88%  '
$ws.Cells.Item(41, 5).Value = '  +1.26%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '14.05'
$ws.Cells.Item(42, 5).Value = '  +4.89%  '
$ws.Cells.Item(43, 2).Value = 'ARBITRUM'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.929'
$ws.Cells.Item(43, 5).Value = '  +1.03%  '
$ws.Cells.Item(44, 2).Value = 'MXToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.74'
$ws.Cells.Item(44, 5).Value = '  +2.33%  '
$ws.Cells.Item(45, 5).Value = '  +2.01%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '6.05'
$ws.Cells.Item(46, 5).Value = '  -0.24%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.07'
$ws.Cells.Item(47, 5).Value = '  +0.10%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.946.06'
$ws.Cells.Item(48, 5).Value = '  -0.01%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '105.57'
$ws.Cells.Item(49, 5).Value = '  -2.02%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0₆0130'
$ws.Cells.Item(50, 5).Value = '  -5.77%  '
$ws.Cells.Item(51, 5).Value = '  -0.03%  '
